$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

# Delete rows 5-7 (the tour now only has 3 matches, rows 2-4)
$ws.Range("A5:AO7").EntireRow.Delete() | Out-Null

# Row 2
$ws.Range("A2").Value = 1369
$ws.Range("B2").Value = '2025-12-19T15:30:00'
$ws.Range("C2").Value = 'Сибирь'
$ws.Range("D2").Value = 'Адмирал'
$ws.Range("E2").Value = 897881
$ws.Range("F2").Value = 'https://text.khl.ru/text/897881.html'
$ws.Range("G2").Value = 1.508387
$ws.Range("H2").Value = 1.2
$ws.Range("I2").Value = 2.714376
$ws.Range("J2").Value = 3.189495
$ws.Range("K2").Value = 2.348941
$ws.Range("L2").Value = 1.957188
$ws.Range("M2").Value = 2.708387
$ws.Range("N2").Value = 23.079923
$ws.Range("O2").Value = 27.688713
$ws.Range("P2").Value = 50.768635
$ws.Range("Q2").Value = -0.123842
$ws.Range("R2").Value = -0.16
$ws.Range("S2").Value = 0.474744
$ws.Range("T2").Value = 0.19579
$ws.Range("U2").Value = 0.329424
$ws.Range("V2").Value = 0.376053
$ws.Range("W2").Value = 0.623905
$ws.Range("X2").Value = 0.569254
$ws.Range("Y2").Value = 0.430704
$ws.Range("Z2").Value = 0.735643
$ws.Range("AA2").Value = 0.264314
$ws.Range("AB2").Value = 0.855059
$ws.Range("AC2").Value = 0.144899
$ws.Range("AD2").Value = 0.928519
$ws.Range("AE2").Value = 0.071439
$ws.Range("AF2").Value = 0.680276
$ws.Range("AG2").Value = 0.319724
$ws.Range("AH2").Value = 0.416896
$ws.Range("AI2").Value = 0.583104
$ws.Range("AJ2").Value = 0.582282
$ws.Range("AK2").Value = 0.417718
$ws.Range("AL2").Value = 0.311737
$ws.Range("AM2").Value = 0.688263
$ws.Range("AN2").Value = 0.826691
$ws.Range("AO2").Value = 0.712628

# Row 3
$ws.Range("A3").Value = 1369
$ws.Range("B3").Value = '2025-12-19T17:00:00'
$ws.Range("C3").Value = 'Трактор'
$ws.Range("D3").Value = 'Металлург Мг'
$ws.Range("E3").Value = 897883
$ws.Range("F3").Value = 'https://text.khl.ru/text/897883.html'
$ws.Range("G3").Value = 5.211111
$ws.Range("H3").Value = 5.288889
$ws.Range("I3").Value = 4.317177
$ws.Range("J3").Value = 2.559277
$ws.Range("K3").Value = 3.885194
$ws.Range("L3").Value = 4.803033
$ws.Range("M3").Value = 10.5
$ws.Range("N3").Value = 35.630589
$ws.Range("O3").Value = 36.669117
$ws.Range("P3").Value = 72.299707
$ws.Range("Q3").Value = 0.16
$ws.Range("R3").Value = 0.16
$ws.Range("S3").Value = 0.311429
$ws.Range("T3").Value = 0.131284
$ws.Range("U3").Value = 0.544567
$ws.Range("V3").Value = 0.026419
$ws.Range("W3").Value = 0.960861
$ws.Range("X3").Value = 0.066438
$ws.Range("Y3").Value = 0.920842
$ws.Range("Z3").Value = 0.135977
$ws.Range("AA3").Value = 0.851303
$ws.Range("AB3").Value = 0.236671
$ws.Range("AC3").Value = 0.750609
$ws.Range("AD3").Value = 0.36165
$ws.Range("AE3").Value = 0.625629
$ws.Range("AF3").Value = 0.899639
$ws.Range("AG3").Value = 0.100361
$ws.Range("AH3").Value = 0.744587
$ws.Range("AI3").Value = 0.255413
$ws.Range("AJ3").Value = 0.952387
$ws.Range("AK3").Value = 0.047613
$ws.Range("AL3").Value = 0.857748
$ws.Range("AM3").Value = 0.142252
$ws.Range("AN3").Value = 0.579935
$ws.Range("AO3").Value = 0.78685

# Row 4
$ws.Range("A4").Value = 1369
$ws.Range("B4").Value = '2025-12-19T19:30:00'
$ws.Range("C4").Value = 'Драконы'
$ws.Range("D4").Value = 'Лада'
$ws.Range("E4").Value = 897882
$ws.Range("F4").Value = 'https://text.khl.ru/text/897882.html'
$ws.Range("G4").Value = 2.820256
$ws.Range("H4").Value = 1.389474
$ws.Range("I4").Value = 3.9
$ws.Range("J4").Value = 3.916575
$ws.Range("K4").Value = 3.368416
$ws.Range("L4").Value = 2.644737
$ws.Range("M4").Value = 4.20973
$ws.Range("N4").Value = 31.730456
$ws.Range("O4").Value = 25.059845
$ws.Range("P4").Value = 56.790301
$ws.Range("Q4").Value = -0.014794
$ws.Range("R4").Value = -0.16
$ws.Range("S4").Value = 0.533076
$ws.Range("T4").Value = 0.159959
$ws.Range("U4").Value = 0.306113
$ws.Range("V4").Value = 0.150034
$ws.Range("W4").Value = 0.849113
$ws.Range("X4").Value = 0.2833
$ws.Range("Y4").Value = 0.715848
$ws.Range("Z4").Value = 0.443569
$ws.Range("AA4").Value = 0.555578
$ws.Range("AB4").Value = 0.60419
$ws.Range("AC4").Value = 0.394957
$ws.Range("AD4").Value = 0.742167
$ws.Range("AE4").Value = 0.256981
$ws.Range("AF4").Value = 0.849534
$ws.Range("AG4").Value = 0.150466
$ws.Range("AH4").Value = 0.654128
$ws.Range("AI4").Value = 0.345872
$ws.Range("AJ4").Value = 0.741136
$ws.Range("AK4").Value = 0.258864
$ws.Range("AL4").Value = 0.492743
$ws.Range("AM4").Value = 0.507257
$ws.Range("AN4").Value = 0.822283
$ws.Range("AO4").Value = 0.630687
